$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.077.44"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.578.01"
$ws.Range("E3").Value = "  +5.96%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'303.75"
$ws.Range("E5").Value = "  +2.81%  "
$ws.Range("D6").Value = "'98.69"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("D7").Value = "'0.596"
$ws.Range("E7").Value = "  +5.26%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.572"
$ws.Range("E9").Value = "  +13.63%  "
$ws.Range("D10").Value = "'38.66"
$ws.Range("E10").Value = "  +11.24%  "
$ws.Range("D11").Value = "'54.29"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").Value = "'0.0835"
$ws.Range("E12").Value = "  +6.70%  "
$ws.Range("D13").Value = "'8.10"
$ws.Range("D14").Value = "2.978.69"
$ws.Range("E14").Value = "  +6.24%  "
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "2.580.25"
$ws.Range("E16").Value = "  +5.58%  "
$ws.Range("D17").Value = "'0.907"
$ws.Range("E17").Value = "  +7.75%  "
$ws.Range("D18").Value = "'14.76"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("D19").Value = "46.282.32"
$ws.Range("E19").Value = "  +1.91%  "
$ws.Range("D20").Value = "0.0₃0999"
$ws.Range("E20").Value = "  +5.97%  "
$ws.Range("D21").Value = "'12.80"
$ws.Range("E21").Value = "  +3.42%  "
$ws.Range("D22").Value = "'6.58"
$ws.Range("E22").Value = "  +6.32%  "
$ws.Range("D23").Value = "'71.32"
$ws.Range("E23").Value = "  +5.97%  "
$ws.Range("D24").Value = "'269.81"
$ws.Range("E24").Value = "  +11.87%  "
$ws.Range("D25").Value = "'2.99"
$ws.Range("E25").Value = "  +7.00%  "
$ws.Range("D26").Value = "'29.64"
$ws.Range("E26").Value = "  +39.18%  "
$ws.Range("E27").Value = "  +9.84%  "
$ws.Range("D28").Value = "'0.996"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "'4.01"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("D30").Value = "'10.43"
$ws.Range("E30").Value = "  +7.75%  "
$ws.Range("D31").Value = "'2.29"
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("D32").Value = "'38.80"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("D33").Value = "'6.09"
$ws.Range("E33").Value = "  +10.76%  "
$ws.Range("D34").Value = "'3.59"
$ws.Range("E34").Value = "  -6.17%  "
$ws.Range("E35").Value = "  +4.10%  "
$ws.Range("D36").Value = "'0.0830"
$ws.Range("E36").Value = "  +8.26%  "
$ws.Range("D37").Value = "'2.15"
$ws.Range("E37").Value = "  +8.03%  "
$ws.Range("D38").Value = "'149.08"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").Value = "'0.120"
$ws.Range("E39").Value = "  +5.76%  "
$ws.Range("D40").Value = "'0.121"
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("D41").Value = "'22.48"
$ws.Range("E41").Value = "  +38.55%  "
$ws.Range("D42").Value = "'15.68"
$ws.Range("E42").Value = "  +6.85%  "
$ws.Range("D43").Value = "'0.0325"
$ws.Range("E43").Value = "  +9.15%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'4.06"
$ws.Range("E44").Value = "  +7.44%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.55"
$ws.Range("E45").Value = "  +9.99%  "
$ws.Range("D46").Value = "2.129.60"
$ws.Range("E46").Value = "  +5.90%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'92.79"
$ws.Range("E48").Value = "  +4.21%  "
$ws.Range("D49").Value = "'9.52"
$ws.Range("E49").Value = "  +10.56%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.77"
$ws.Range("E50").Value = "  +3.04%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'107.73"
$ws.Range("E51").Value = "  +6.47%  "
